$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "平潭发展"
$ws.Range("C2").Value = "药明康德"
$ws.Range("A3").Value = "方正科技"
$ws.Range("B3").Value = "平潭发展"
$ws.Range("C3").Value = "华工科技"
$ws.Range("A4").Value = "中国卫星"
$ws.Range("B4").Value = "中国卫星"
$ws.Range("C4").Value = "剑桥科技"
$ws.Range("A5").Value = "香农芯创"
$ws.Range("B5").Value = "东方财富"
$ws.Range("A6").Value = "达华智能"
$ws.Range("B6").Value = "达华智能"
$ws.Range("C6").Value = "三花智控"
$ws.Range("A7").Value = "药明康德"
$ws.Range("B7").Value = "闻泰科技"
$ws.Range("C7").Value = "平潭发展"
$ws.Range("A8").Value = "闻泰科技"
$ws.Range("B8").Value = "香农芯创"
$ws.Range("C8").Value = "盈新发展"
$ws.Range("A9").Value = "盈新发展"
$ws.Range("B9").Value = "药明康德"
$ws.Range("C9").Value = "闻泰科技"
$ws.Range("A10").Value = "大洋电机"
$ws.Range("C10").Value = "大洋电机"
$ws.Range("A11").Value = "航天科技"
$ws.Range("B11").Value = "大为股份"
$ws.Range("C11").Value = "达华智能"
$ws.Range("A12").Value = "晶瑞电材"
$ws.Range("B12").Value = "中国核建"
$ws.Range("C12").Value = "超颖电子"
$ws.Range("A13").Value = "中际旭创"
$ws.Range("B13").Value = "晶瑞电材"
$ws.Range("A14").Value = "三花智控"
$ws.Range("B14").Value = "格尔软件"
$ws.Range("C14").Value = "英唐智控"
$ws.Range("A15").Value = "大为股份"
$ws.Range("B15").Value = "大洋电机"
$ws.Range("C15").Value = "神开股份"
$ws.Range("A16").Value = "和而泰"
$ws.Range("B16").Value = "上海电力"
$ws.Range("C16").Value = "国盛金控"
$ws.Range("A17").Value = "中国核建"
$ws.Range("B17").Value = "航天科技"
$ws.Range("C17").Value = "山子高科"
$ws.Range("A18").Value = "东方财富"
$ws.Range("B18").Value = "和而泰"
$ws.Range("C18").Value = "上海电力"
$ws.Range("A19").Value = "上海电力"
$ws.Range("B19").Value = "三花智控"
$ws.Range("C19").Value = "大有能源"
$ws.Range("A20").Value = "江波龙"
$ws.Range("B20").Value = "深科技"
$ws.Range("C20").Value = "中际旭创"
$ws.Range("A21").Value = "贵州茅台"
$ws.Range("B21").Value = "神州信息"
$ws.Range("C21").Value = "大为股份"
